$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (Strike#) values. Regenerate them (the new K values
# replace the old Strike# values) for rows 2-8.
$kValues = @{
    2 = 4
    3 = 1
    4 = 1
    5 = 4
    6 = 1
    7 = 1
    8 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
